# Update the "想去人数" (want-to-go count) figures in the F column for both
# the "展览" sheet and the "全部类型" sheet, which mirror the same data.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 1886
    "F3"  = 82
    "F5"  = 196
    "F6"  = 762
    "F9"  = 4518
    "F12" = 1291
    "F15" = 878
    "F17" = 489
    "F19" = 233
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
